$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Temperatura (D) and Nori (E) values for rows 2-6
$ws.Range("D2").Value = 5.43
$ws.Range("E2").Value = 71

$ws.Range("D3").Value = 5.44
$ws.Range("E3").Value = 72

$ws.Range("D4").Value = 5.55
$ws.Range("E4").Value = 75

$ws.Range("D5").Value = 5.89
$ws.Range("E5").Value = 88

$ws.Range("D6").Value = 5.96
$ws.Range("E6").Value = 94

# Populate Radiatie (C) values for rows 88-97 (previously blank)
$ws.Range("C88").Value = 0
$ws.Range("C89").Value = 0
$ws.Range("C90").Value = 0
$ws.Range("C91").Value = 0
$ws.Range("C92").Value = 0
$ws.Range("C93").Value = 0
$ws.Range("C94").Value = 15.77
$ws.Range("C95").Value = 133.52
$ws.Range("C96").Value = 305.36
$ws.Range("C97").Value = 476.76
